# Update "想去人数" (number of people interested) figures to the latest
# scraped values, as produced by the gh-pages data generation run at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 14155
$wsExhibit.Range("F3").Value = 560
$wsExhibit.Range("F6").Value = 1053
$wsExhibit.Range("F7").Value = 13957
$wsExhibit.Range("F8").Value = 15090
$wsExhibit.Range("F10").Value = 22
$wsExhibit.Range("F14").Value = 47
$wsExhibit.Range("F20").Value = 23
$wsExhibit.Range("F21").Value = 1170
$wsExhibit.Range("F24").Value = 5863
$wsExhibit.Range("F27").Value = 5472
$wsExhibit.Range("F30").Value = 72
$wsExhibit.Range("F31").Value = 361

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 10

# --- Sheet "全部类型" (All types, aggregated) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14155
$wsAll.Range("F3").Value = 10
$wsAll.Range("F4").Value = 560
$wsAll.Range("F7").Value = 1053
$wsAll.Range("F8").Value = 13957
$wsAll.Range("F9").Value = 15090
$wsAll.Range("F11").Value = 22
$wsAll.Range("F15").Value = 47
$wsAll.Range("F21").Value = 23
$wsAll.Range("F22").Value = 1170
$wsAll.Range("F26").Value = 5863
$wsAll.Range("F29").Value = 5472
$wsAll.Range("F32").Value = 72
$wsAll.Range("F33").Value = 361
